# "Generate Report for Handoff"
# The localization status report is regenerated: the in-progress status
# ("In Translation") becomes "Ready for handoff" everywhere it is shown,
# and the handoff-generation timestamps are bumped to the new run time.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$ws_overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$ws_overview.Range("F2").Value = "Ready for handoff"   # de-de status
$ws_zhcn.Range("C2").Value     = "Ready for handoff"   # Status column
$ws_dede.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime ---------------
$ws_overview.Range("G2").Value = "2016-09-01 22:43:10" # also de-de's handoff datetime
$ws_dede.Range("H2").Value     = "2016-09-01 22:43:10"
$ws_zhcn.Range("H2").Value     = "2016-09-01 22:42:59"

# --- Column widths: the "Status" columns grew to fit "Ready for handoff" ---
$ws_overview.Columns.Item(5).ColumnWidth = 16.26
$ws_overview.Columns.Item(6).ColumnWidth = 16.26
$ws_zhcn.Columns.Item(3).ColumnWidth     = 16.26
$ws_dede.Columns.Item(3).ColumnWidth     = 16.26
